$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: update the "Action recommended" text for the heading issue ---
$ws.Range("E8").Value = "adding a h2 tag on Let's talk web design!"

# --- Row 13 (new): Contact Me page .min links/scripts issue ---
# Match formatting of the neighboring rows first (Copy() also copies the source value/text)
$ws.Range("B9").Copy($ws.Range("B13"))
$ws.Range("C9").Copy($ws.Range("C13"))
$ws.Range("E12").Copy($ws.Range("E13"))
$ws.Rows.Item(13).RowHeight = 16.5

# Now set the real content for the new row
$ws.Range("A13").Value = "Accessibility"
$ws.Range("B13").Value = "some of the links and scripts in Contact Me page is not working like bootstrap, javascript and font awesome that makes the site not functional "
$ws.Range("C13").Value = "The links and script in Contact Me page has .min version in src and href that in files we don’t have the minimized version of the code"
$ws.Range("E13").Value = "I recommend to delete all .min from all the links that has this word in order to make the codes apply to the page  since we have the codes that are not minimized"
$ws.Range("F13").Value = "We don't need a Reference for this problem since its just a typing mistake and there is no minimized version of the code"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 108.85714285714286
$ws.Columns.Item(3).ColumnWidth = 103.42857142857143
$ws.Columns.Item(5).ColumnWidth = 116.0
$ws.Columns.Item(6).ColumnWidth = 92.14285714285714

# --- Sheet view: zoom + selection ---
$ws.Range("B14").Select()
$excel.ActiveWindow.Zoom = 84
